$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-CellText $ws.Range('D2') '26.915.86'
Set-CellText $ws.Range('E2') '  -1.12%  '

Set-CellText $ws.Range('D3') '1.738.09'
Set-CellText $ws.Range('E3') '  +0.98%  '

Set-CellText $ws.Range('D4') '1.001'
Set-CellText $ws.Range('E4') '  -0.02%  '

Set-CellText $ws.Range('D5') '311.22'
Set-CellText $ws.Range('E5') '  -0.74%  '

Set-CellText $ws.Range('E6') '  -0.01%  '

Set-CellText $ws.Range('D7') '0.5051'
Set-CellText $ws.Range('E7') '  +9.78%  '

Set-CellText $ws.Range('D8') '0.3560'
Set-CellText $ws.Range('E8') '  +3.44%  '

Set-CellText $ws.Range('D9') '42.09'
Set-CellText $ws.Range('E9') '  -0.73%  '

Set-CellText $ws.Range('D10') '0.07249'

Set-CellText $ws.Range('D11') '1.060'
Set-CellText $ws.Range('E11') '  +1.45%  '

Set-CellText $ws.Range('E12') '  +0.10%  '

Set-CellText $ws.Range('D13') '20.22'
Set-CellText $ws.Range('E13') '  +2.29%  '

Set-CellText $ws.Range('D14') '5.939'
Set-CellText $ws.Range('E14') '  +1.86%  '

Set-CellText $ws.Range('D15') '1.741.05'
Set-CellText $ws.Range('E15') '  +1.07%  '

Set-CellText $ws.Range('D16') '6.820'
Set-CellText $ws.Range('E16') '  -0.68%  '

Set-CellText $ws.Range('D17') '86.76'
Set-CellText $ws.Range('E17') '  -2.49%  '

Set-CellText $ws.Range('D18') '0.00001034'
Set-CellText $ws.Range('E18') '  -0.51%  '

Set-CellText $ws.Range('D19') '0.06408'
Set-CellText $ws.Range('E19') '  +1.22%  '

Set-CellText $ws.Range('E20') '  +0.04%  '

Set-CellText $ws.Range('D21') '16.47'
Set-CellText $ws.Range('E21') '  -0.13%  '

Set-CellText $ws.Range('D22') '5.711'
Set-CellText $ws.Range('E22') '  +1.88%  '

Set-CellText $ws.Range('D23') '26.984.01'
Set-CellText $ws.Range('E23') '  -0.97%  '

Set-CellText $ws.Range('E24') '  +4.38%  '

Set-CellText $ws.Range('D25') '2.043'
Set-CellText $ws.Range('E25') '  -4.17%  '

Set-CellText $ws.Range('D26') '154.03'
Set-CellText $ws.Range('E26') '  -0.25%  '

Set-CellText $ws.Range('D27') '19.87'
Set-CellText $ws.Range('E27') '  +2.84%  '

Set-CellText $ws.Range('D28') '1.942.02'
Set-CellText $ws.Range('E28') '  +1.12%  '

Set-CellText $ws.Range('D29') '2.206'
Set-CellText $ws.Range('E29') '  +2.95%  '

Set-CellText $ws.Range('D30') '120.01'
Set-CellText $ws.Range('E30') '  +0.61%  '

Set-CellText $ws.Range('D31') '1.044'
Set-CellText $ws.Range('E31') '  +1.54%  '

Set-CellText $ws.Range('D32') '0.09505'
Set-CellText $ws.Range('E32') '  +4.60%  '

Set-CellText $ws.Range('D33') '3.581'
Set-CellText $ws.Range('E33') '  -0.50%  '

Set-CellText $ws.Range('D34') '5.367'
Set-CellText $ws.Range('E34') '  +0.54%  '

Set-CellText $ws.Range('D35') '0.02195'
Set-CellText $ws.Range('E35') '  -0.57%  '

Set-CellText $ws.Range('D36') '0.05824'
Set-CellText $ws.Range('E36') '  -0.21%  '

Set-CellText $ws.Range('D37') '11.08'
Set-CellText $ws.Range('E37') '  +0.31%  '

Set-CellText $ws.Range('B38') 'Algorand'
Set-CellText $ws.Range('C38') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-CellText $ws.Range('D38') '0.1999'
Set-CellText $ws.Range('E38') '  +0.13%  '

Set-CellText $ws.Range('B39') 'WEMIXTOKEN'
Set-CellText $ws.Range('C39') 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-CellText $ws.Range('D39') '1.424'
Set-CellText $ws.Range('E39') '  +1.05%  '

Set-CellText $ws.Range('D40') '4.776'
Set-CellText $ws.Range('E40') '  +0.77%  '

Set-CellText $ws.Range('D41') '0.6040'
Set-CellText $ws.Range('E41') '  +1.76%  '

Set-CellText $ws.Range('D42') '1.108'
Set-CellText $ws.Range('E42') '  -1.72%  '

Set-CellText $ws.Range('D43') '7.584'
Set-CellText $ws.Range('E43') '  +1.64%  '

Set-CellText $ws.Range('D44') '12.76'
Set-CellText $ws.Range('E44') '  -0.07%  '

Set-CellText $ws.Range('D45') '3.595'
Set-CellText $ws.Range('E45') '  -0.15%  '

Set-CellText $ws.Range('D46') '0.5650'
Set-CellText $ws.Range('E46') '  +0.12%  '

Set-CellText $ws.Range('D47') '120.24'
Set-CellText $ws.Range('E47') '  +0.55%  '

Set-CellText $ws.Range('D48') '1.853'
Set-CellText $ws.Range('E48') '  -0.66%  '

Set-CellText $ws.Range('B49') 'EOS'
Set-CellText $ws.Range('C49') 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-CellText $ws.Range('D49') '1.102'
Set-CellText $ws.Range('E49') '  +1.82%  '

Set-CellText $ws.Range('B50') 'Cronos'
Set-CellText $ws.Range('C50') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText $ws.Range('D50') '0.06666'
Set-CellText $ws.Range('E50') '  +0.15%  '

Set-CellText $ws.Range('E51') '  +0.03%  '
